$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23
$ws.Cells.Item(23, 2).Value = 6893102
$ws.Cells.Item(23, 6).Value = 'UD Oliveirense'
$ws.Cells.Item(23, 7).Value = 'Penafiel'
$ws.Cells.Item(23, 8).Value = 3
$ws.Cells.Item(23, 11).Value = 2.25
$ws.Cells.Item(23, 13).Value = 2.75
$ws.Cells.Item(23, 14).Value = 2.25
$ws.Cells.Item(23, 15).Value = 3.3
$ws.Cells.Item(23, 16).Value = 2.75
$ws.Cells.Item(23, 17).Value = -0.25
$ws.Cells.Item(23, 18).Value = 2.025
$ws.Cells.Item(23, 19).Value = 1.775
$ws.Cells.Item(23, 22).Value = 1.775
$ws.Cells.Item(23, 23).Value = 1.25
$ws.Cells.Item(23, 26).Value = 1.025
# Row 24
$ws.Cells.Item(24, 2).Value = 6899210
$ws.Cells.Item(24, 6).Value = 'FC Porto B'
$ws.Cells.Item(24, 7).Value = 'UD Leiria'
$ws.Cells.Item(24, 8).Value = 2
$ws.Cells.Item(24, 11).Value = 2.1
$ws.Cells.Item(24, 13).Value = 3
$ws.Cells.Item(24, 14).Value = 2.45
$ws.Cells.Item(24, 15).Value = 3.2
$ws.Cells.Item(24, 16).Value = 2.5
$ws.Cells.Item(24, 17).Value = 0
$ws.Cells.Item(24, 18).Value = 1.9
$ws.Cells.Item(24, 19).Value = 1.95
$ws.Cells.Item(24, 22).Value = 1.825
$ws.Cells.Item(24, 23).Value = 1.45
$ws.Cells.Item(24, 26).Value = 0.8999999999999999
# Row 39
$ws.Cells.Item(39, 2).Value = 6899293
$ws.Cells.Item(39, 6).Value = 'AVS'
$ws.Cells.Item(39, 7).Value = 'Vilaverdense'
$ws.Cells.Item(39, 8).Value = 2
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 'H'
$ws.Cells.Item(39, 11).Value = 1.8
$ws.Cells.Item(39, 12).Value = 3.25
$ws.Cells.Item(39, 13).Value = 4
$ws.Cells.Item(39, 14).Value = 1.727
$ws.Cells.Item(39, 15).Value = 3.4
$ws.Cells.Item(39, 16).Value = 4.2
$ws.Cells.Item(39, 17).Value = -0.75
$ws.Cells.Item(39, 19).Value = 1.8
$ws.Cells.Item(39, 20).Value = 2.5
$ws.Cells.Item(39, 21).Value = 2
$ws.Cells.Item(39, 22).Value = 1.85
$ws.Cells.Item(39, 23).Value = 0.7270000000000001
$ws.Cells.Item(39, 25).Value = -1
$ws.Cells.Item(39, 26).Value = 1.05
$ws.Cells.Item(39, 27).Value = -1
$ws.Cells.Item(39, 28).Value = -1
$ws.Cells.Item(39, 29).Value = 0.8500000000000001
# Row 40
$ws.Cells.Item(40, 2).Value = 6899208
$ws.Cells.Item(40, 6).Value = 'FC Porto B'
$ws.Cells.Item(40, 7).Value = 'Maritimo'
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 2
$ws.Cells.Item(40, 10).Value = 'A'
$ws.Cells.Item(40, 11).Value = 2.55
$ws.Cells.Item(40, 12).Value = 3
$ws.Cells.Item(40, 13).Value = 2.55
$ws.Cells.Item(40, 14).Value = 2.3
$ws.Cells.Item(40, 15).Value = 3
$ws.Cells.Item(40, 16).Value = 2.9
$ws.Cells.Item(40, 17).Value = -0.25
$ws.Cells.Item(40, 19).Value = 1.75
$ws.Cells.Item(40, 20).Value = 2.25
$ws.Cells.Item(40, 21).Value = 1.775
$ws.Cells.Item(40, 22).Value = 2.025
$ws.Cells.Item(40, 23).Value = -1
$ws.Cells.Item(40, 25).Value = 1.9
$ws.Cells.Item(40, 26).Value = -1
$ws.Cells.Item(40, 27).Value = 0.75
$ws.Cells.Item(40, 28).Value = -0.5
$ws.Cells.Item(40, 29).Value = 0.5125
# Row 186
$ws.Cells.Item(186, 2).Value = 6893598
$ws.Cells.Item(186, 6).Value = 'FC Porto B'
$ws.Cells.Item(186, 7).Value = 'CD Mafra'
$ws.Cells.Item(186, 8).Value = 1
$ws.Cells.Item(186, 9).Value = 1
$ws.Cells.Item(186, 11).Value = 2.05
$ws.Cells.Item(186, 12).Value = 3.4
$ws.Cells.Item(186, 13).Value = 3.4
$ws.Cells.Item(186, 14).Value = 1.909
$ws.Cells.Item(186, 15).Value = 3.5
$ws.Cells.Item(186, 16).Value = 3.75
$ws.Cells.Item(186, 17).Value = -0.5
$ws.Cells.Item(186, 18).Value = 1.9
$ws.Cells.Item(186, 19).Value = 1.9
$ws.Cells.Item(186, 20).Value = 2.5
$ws.Cells.Item(186, 21).Value = 1.85
$ws.Cells.Item(186, 22).Value = 1.95
$ws.Cells.Item(186, 24).Value = 2.5
$ws.Cells.Item(186, 26).Value = -1
$ws.Cells.Item(186, 27).Value = 0.8999999999999999
$ws.Cells.Item(186, 29).Value = 0.95
# Row 187
$ws.Cells.Item(187, 2).Value = 6893600
$ws.Cells.Item(187, 6).Value = 'CF Os Belenenses'
$ws.Cells.Item(187, 7).Value = 'Tondela'
$ws.Cells.Item(187, 8).Value = 0
$ws.Cells.Item(187, 9).Value = 0
$ws.Cells.Item(187, 11).Value = 3.8
$ws.Cells.Item(187, 12).Value = 3.3
$ws.Cells.Item(187, 13).Value = 1.95
$ws.Cells.Item(187, 14).Value = 3.5
$ws.Cells.Item(187, 15).Value = 3.3
$ws.Cells.Item(187, 16).Value = 2.05
$ws.Cells.Item(187, 17).Value = 0.25
$ws.Cells.Item(187, 18).Value = 2.05
$ws.Cells.Item(187, 19).Value = 1.8
$ws.Cells.Item(187, 20).Value = 2.25
$ws.Cells.Item(187, 21).Value = 1.95
$ws.Cells.Item(187, 22).Value = 1.9
$ws.Cells.Item(187, 24).Value = 2.3
$ws.Cells.Item(187, 26).Value = 0.5249999999999999
$ws.Cells.Item(187, 27).Value = -0.5
$ws.Cells.Item(187, 29).Value = 0.8999999999999999
# Row 230
$ws.Cells.Item(230, 2).Value = 6893180
$ws.Cells.Item(230, 6).Value = 'Academico Viseu'
$ws.Cells.Item(230, 7).Value = 'Pacos Ferreira'
$ws.Cells.Item(230, 8).Value = 1
$ws.Cells.Item(230, 9).Value = 1
$ws.Cells.Item(230, 10).Value = 'D'
$ws.Cells.Item(230, 11).Value = 2.375
$ws.Cells.Item(230, 12).Value = 3
$ws.Cells.Item(230, 13).Value = 2.8
$ws.Cells.Item(230, 14).Value = 2.7
$ws.Cells.Item(230, 15).Value = 2.9
$ws.Cells.Item(230, 16).Value = 2.5
$ws.Cells.Item(230, 17).Value = 0
$ws.Cells.Item(230, 21).Value = 1.975
$ws.Cells.Item(230, 22).Value = 1.825
$ws.Cells.Item(230, 23).Value = -1
$ws.Cells.Item(230, 24).Value = 1.9
$ws.Cells.Item(230, 26).Value = 0
$ws.Cells.Item(230, 27).Value = 0
$ws.Cells.Item(230, 29).Value = 0.4125
# Row 231
$ws.Cells.Item(231, 2).Value = 6899162
$ws.Cells.Item(231, 6).Value = 'Maritimo'
$ws.Cells.Item(231, 7).Value = 'UD Leiria'
$ws.Cells.Item(231, 8).Value = 2
$ws.Cells.Item(231, 9).Value = 0
$ws.Cells.Item(231, 10).Value = 'H'
$ws.Cells.Item(231, 11).Value = 1.8
$ws.Cells.Item(231, 12).Value = 3.25
$ws.Cells.Item(231, 13).Value = 4
$ws.Cells.Item(231, 14).Value = 1.727
$ws.Cells.Item(231, 15).Value = 3.3
$ws.Cells.Item(231, 16).Value = 4.333
$ws.Cells.Item(231, 17).Value = -0.75
$ws.Cells.Item(231, 21).Value = 1.8
$ws.Cells.Item(231, 22).Value = 2
$ws.Cells.Item(231, 23).Value = 0.7270000000000001
$ws.Cells.Item(231, 24).Value = -1
$ws.Cells.Item(231, 26).Value = 1
$ws.Cells.Item(231, 27).Value = -1
$ws.Cells.Item(231, 29).Value = 0.5
# Row 236
$ws.Cells.Item(236, 21).Value = 1.975
$ws.Cells.Item(236, 22).Value = 1.875
# Row 237
$ws.Cells.Item(237, 21).Value = 1.85
$ws.Cells.Item(237, 22).Value = 2
# Row 238
$ws.Cells.Item(238, 18).Value = 1.8
$ws.Cells.Item(238, 19).Value = 2.05
